$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 534, shifting existing rows 534:586 down to 535:587
$ws.Rows.Item(534).Insert()

# Populate the newly inserted row 534 with the new data record
$ws.Cells.Item(534, 1).Value = 11
$ws.Cells.Item(534, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(534, 3).Value = "Bíobío"
$ws.Cells.Item(534, 4).Value = 45212
$ws.Cells.Item(534, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(534, 5).Value = 8
$ws.Cells.Item(534, 6).Value = 100114014
$ws.Cells.Item(534, 7).Value = "Betarraga"
$ws.Cells.Item(534, 8).Value = "Sin especificar"
$ws.Cells.Item(534, 9).Value = "Primera"
$ws.Cells.Item(534, 10).Value = 650
$ws.Cells.Item(534, 11).Value = 600
$ws.Cells.Item(534, 12).Value = 650
$ws.Cells.Item(534, 13).Value = 627
$ws.Cells.Item(534, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(534, 15).Value = "Región Metropolitana"
$ws.Cells.Item(534, 16).Value = 125
$ws.Cells.Item(534, 17).Value = 5
$ws.Cells.Item(534, 18).Value = "Hortaliza"
